$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($cell, $val)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextCell $ws.Range("B23") "Dai"
Set-TextCell $ws.Range("C23") "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextCell $ws.Range("D23") "0.998"
Set-TextCell $ws.Range("E23") "  -0.16%  "

Set-TextCell $ws.Range("B24") "NEARProtocol"
Set-TextCell $ws.Range("C24") "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextCell $ws.Range("D24") "4.34"
Set-TextCell $ws.Range("E24") "  -7.90%  "

Set-TextCell $ws.Range("B42") "PolygonEcosystemToken"
Set-TextCell $ws.Range("C42") "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextCell $ws.Range("D42") "0.331"
Set-TextCell $ws.Range("E42") "  -7.21%  "

Set-TextCell $ws.Range("B43") "USDe"
Set-TextCell $ws.Range("C43") "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextCell $ws.Range("D43") "1.00"
Set-TextCell $ws.Range("E43") "  -0.04%  "

Set-TextCell $ws.Range("B44") "RenderToken"
Set-TextCell $ws.Range("C44") "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
Set-TextCell $ws.Range("D44") "4.83"
Set-TextCell $ws.Range("E44") "  -7.27%  "

Set-TextCell $ws.Range("D2") "66.385.46"
Set-TextCell $ws.Range("E2") "  -2.64%  "
Set-TextCell $ws.Range("D3") "2.447.54"
Set-TextCell $ws.Range("E3") "  -3.31%  "
Set-TextCell $ws.Range("E4") "  +0.12%  "
Set-TextCell $ws.Range("D5") "578.14"
Set-TextCell $ws.Range("E5") "  -2.60%  "
Set-TextCell $ws.Range("D6") "162.69"
Set-TextCell $ws.Range("E6") "  -8.69%  "
Set-TextCell $ws.Range("E7") "  +0.19%  "
Set-TextCell $ws.Range("D8") "0.510"
Set-TextCell $ws.Range("E8") "  -4.06%  "
Set-TextCell $ws.Range("D9") "2.456.49"
Set-TextCell $ws.Range("E9") "  -2.97%  "
Set-TextCell $ws.Range("D10") "0.134"
Set-TextCell $ws.Range("E10") "  -6.14%  "
Set-TextCell $ws.Range("D11") "0.164"
Set-TextCell $ws.Range("E11") "  -0.24%  "
Set-TextCell $ws.Range("D12") "0.336"
Set-TextCell $ws.Range("E12") "  -2.89%  "
Set-TextCell $ws.Range("D13") "4.89"
Set-TextCell $ws.Range("E13") "  -4.54%  "
Set-TextCell $ws.Range("D14") "25.52"
Set-TextCell $ws.Range("E14") "  -4.93%  "
Set-TextCell $ws.Range("D15") "2.911.53"
Set-TextCell $ws.Range("E15") "  -2.75%  "
Set-TextCell $ws.Range("D16") "0.0000171"
Set-TextCell $ws.Range("E16") "  -5.36%  "
Set-TextCell $ws.Range("D17") "66.501.19"
Set-TextCell $ws.Range("E17") "  -2.49%  "
Set-TextCell $ws.Range("D18") "2.469.08"
Set-TextCell $ws.Range("E18") "  -3.35%  "
Set-TextCell $ws.Range("D19") "11.41"
Set-TextCell $ws.Range("E19") "  -1.11%  "
Set-TextCell $ws.Range("D20") "7.69"
Set-TextCell $ws.Range("E20") "  -3.79%  "
Set-TextCell $ws.Range("D21") "354.52"
Set-TextCell $ws.Range("E21") "  -3.54%  "
Set-TextCell $ws.Range("D22") "4.03"
Set-TextCell $ws.Range("E22") "  -4.19%  "
Set-TextCell $ws.Range("D25") "70.21"
Set-TextCell $ws.Range("E25") "  -0.85%  "
Set-TextCell $ws.Range("D26") "1.82"
Set-TextCell $ws.Range("E26") "  -5.82%  "
Set-TextCell $ws.Range("D27") "9.24"
Set-TextCell $ws.Range("E27") "  -9.64%  "
Set-TextCell $ws.Range("D28") "0.998"
Set-TextCell $ws.Range("E28") "  -0.07%  "
Set-TextCell $ws.Range("D29") "2.606.08"
Set-TextCell $ws.Range("E29") "  -2.44%  "
Set-TextCell $ws.Range("D30") "0.0₃0916"
Set-TextCell $ws.Range("E30") "  -8.05%  "
Set-TextCell $ws.Range("D31") "7.94"
Set-TextCell $ws.Range("E31") "  -4.24%  "
Set-TextCell $ws.Range("D32") "493.07"
Set-TextCell $ws.Range("E32") "  -8.94%  "
Set-TextCell $ws.Range("D33") "1.81"
Set-TextCell $ws.Range("E33") "  -3.54%  "
Set-TextCell $ws.Range("D34") "1.25"
Set-TextCell $ws.Range("E34") "  -6.76%  "
Set-TextCell $ws.Range("E35") "  +0.24%  "
Set-TextCell $ws.Range("D36") "0.125"
Set-TextCell $ws.Range("E36") "  -3.52%  "
Set-TextCell $ws.Range("D37") "159.16"
Set-TextCell $ws.Range("E37") "  +1.32%  "
Set-TextCell $ws.Range("D38") "18.62"
Set-TextCell $ws.Range("E38") "  -1.26%  "
Set-TextCell $ws.Range("D39") "1.39"
Set-TextCell $ws.Range("E39") "  -4.99%  "
Set-TextCell $ws.Range("D40") "18.53"
Set-TextCell $ws.Range("E40") "  -0.83%  "
Set-TextCell $ws.Range("D41") "1.70"
Set-TextCell $ws.Range("E41") "  -6.36%  "
Set-TextCell $ws.Range("D45") "2.41"
Set-TextCell $ws.Range("E45") "  -6.11%  "
Set-TextCell $ws.Range("D46") "39.19"
Set-TextCell $ws.Range("E46") "  -2.11%  "
Set-TextCell $ws.Range("D47") "139.80"
Set-TextCell $ws.Range("E47") "  -5.18%  "
Set-TextCell $ws.Range("D48") "3.55"
Set-TextCell $ws.Range("E48") "  -4.65%  "
Set-TextCell $ws.Range("D49") "0.527"
Set-TextCell $ws.Range("E49") "  -6.08%  "
Set-TextCell $ws.Range("D50") "0.0₆0257"
Set-TextCell $ws.Range("E50") "  -7.37%  "
Set-TextCell $ws.Range("D51") "1.62"
Set-TextCell $ws.Range("E51") "  -5.34%  "
